$wb = $excel.ActiveWorkbook

# =========================================================================
# Sheet "展览" (Exhibitions) - update "想去人数" (want-to-go count) values
# =========================================================================
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 6064
$ws.Range("F10").Value = 701
$ws.Range("F11").Value = 1748
$ws.Range("F12").Value = 1748
$ws.Range("F13").Value = 8
$ws.Range("F14").Value = 1682
$ws.Range("F15").Value = 584
$ws.Range("F16").Value = 214
$ws.Range("F17").Value = 662
$ws.Range("F18").Value = 4784
$ws.Range("F19").Value = 128
$ws.Range("F21").Value = 682
$ws.Range("F24").Value = 27
$ws.Range("F25").Value = 58
$ws.Range("F26").Value = 27
$ws.Range("F27").Value = 2371
$ws.Range("F28").Value = 53
$ws.Range("F30").Value = 14
$ws.Range("F33").Value = 1253
$ws.Range("F34").Value = 793
$ws.Range("F35").Value = 35
$ws.Range("F36").Value = 12
$ws.Range("F38").Value = 1324
$ws.Range("F39").Value = 1306
$ws.Range("F40").Value = 88

# =========================================================================
# Sheet "演出" (Performances) - update "想去人数" values
# =========================================================================
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 3
$ws.Range("F5").Value = 116
$ws.Range("F11").Value = 17
$ws.Range("F14").Value = 95
$ws.Range("F20").Value = 316
$ws.Range("F21").Value = 245
$ws.Range("F22").Value = 508

# =========================================================================
# Sheet "本地生活" (Local life) - update values + mark row 5 sold out
# =========================================================================
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 797
$ws.Range("F4").Value = 224
$ws.Range("F5").Value = 322
$ws.Range("G5").Value = "已售罄"

# =========================================================================
# Sheet "全部类型" (All types) - update "想去人数" values
# =========================================================================
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 797
$ws.Range("F6").Value = 224

# Rows 7, 8 and 9 are fully refreshed with newer source data (the feed this
# sheet mirrors was re-pulled, rotating which events land on which row and
# adding a previously-missing event in row 9).

# Row 7: now the "ICOS国际动漫节×CGF中国游戏节04" entry
$ws.Range("B7").Value = "'2024-10-19"
$ws.Range("C7").Value = "北京·ICOS国际动漫节×CGF中国游戏节04"
$ws.Range("D7").Value = "石景山路68号 北京首钢会展中心"
$ws.Range("E7").Value = "2024.10.19 09:00-10.20 17:00"
$ws.Range("F7").Value = 6064
$ws.Range("G7").Value = 90
$ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=88085"
$ws.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202406/jQr9LeQO1719381394199.jpeg"

# Row 8: now the "春江花月夜-大师的启蒙古筝专场音乐会" entry
$ws.Range("B8").Value = "'2024-10-20"
$ws.Range("C8").Value = " 北京·《春江花月夜-大师的启蒙古筝专场音乐会》"
$ws.Range("D8").Value = "北京东图剧场 北京东图剧场"
$ws.Range("E8").Value = "2024.10.20 19:30-10.20 21:00"
$ws.Range("F8").Value = 3
$ws.Range("G8").Value = 63
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=92878"
$ws.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202409/tH4Z9ch71727197609455.png"

# Row 9: now the "光辉岁月-大师的启蒙屈红震架子鼓专场演奏会" entry (new row)
$ws.Range("C9").Value = "北京·《光辉岁月-大师的启蒙屈红震架子鼓专场演奏会》"
$ws.Range("E9").Value = "2024.10.20 14:30-10.20 16:00"
$ws.Range("F9").Value = 0
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=92872"
$ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202409/b9A9LSUc1727198181199.png"

$ws.Range("F10").Value = 116
$ws.Range("F18").Value = 17
$ws.Range("F21").Value = 1748
$ws.Range("F22").Value = 8
$ws.Range("F23").Value = 1682
$ws.Range("F24").Value = 96
$ws.Range("F25").Value = 584
$ws.Range("F26").Value = 214
$ws.Range("F27").Value = 662
$ws.Range("F28").Value = 4784
$ws.Range("F30").Value = 682
$ws.Range("F32").Value = 27
$ws.Range("F33").Value = 58
$ws.Range("F35").Value = 27
$ws.Range("F36").Value = 2371
$ws.Range("F38").Value = 14
$ws.Range("F40").Value = 1253
$ws.Range("F42").Value = 245
$ws.Range("F43").Value = 508
$ws.Range("F44").Value = 793
$ws.Range("F45").Value = 35
$ws.Range("F46").Value = 12
$ws.Range("F48").Value = 1324
$ws.Range("F50").Value = 88
